$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 2 is inserted (LINExChatGPTx... bot posting); all previously existing rows
# shift down by one and pick up the refreshed scrape timestamp. Rewriting A2:H20 directly
# (rather than Rows.Insert, whose Hyperlinks bookkeeping this host does not shift) keeps
# data + hyperlinks consistent in one pass.
$rows = @(
    ,@("2025-09-26 01:42:45", "LINExChatGPTx美容室向け予約Bot (仕様書、契約書あり)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400801", 445, "🔥GPT,ChatGPT ★bot")
    ,@("2025-09-26 01:42:45", "【SES案件多数】バックエンドエンジニア募集(Java/PHP/Python/Node.js)", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399874", 320, "🔥Python ★Java ◆Node.js ○PHP")
    ,@("2025-09-26 01:42:45", "【低予算希望】LINE公式アカウント+社食注文システム開発依頼(社内利用のみ)", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400375", 118, "◆開発,システム開発")
    ,@("2025-09-26 01:42:45", "【急募】LLMによるMCP(Model Context Protocol)でのExcel操作機能開発", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400689", 75, "◆開発")
    ,@("2025-09-26 01:42:45", "【急募】住宅展示場マッチング診断サービスのMVP開発依頼", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399759", 75, "◆開発")
    ,@("2025-09-26 01:42:45", "【フリーランス募集】CTビューアーソフト気道抽出機能開発", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400101", 68, "◆開発")
    ,@("2025-09-26 01:42:45", "【急募】音源ライセンス販売サイトのMVP構築依頼", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400763", 45, "◇サイト")
    ,@("2025-09-26 01:42:45", "MYSQLからGoogleスプレッドシートへデータ取り込み及びスプレッドシート改修", "システム開発", "10,000 円 ~ 20,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400606", 30, "◇MySQL")
    ,@("2025-09-26 01:42:45", "eBayテラピークでのキーワード検索結果等の取得するためのシステム制作", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5390238", 33, $null)
    ,@("2025-09-26 01:42:45", "Drupal関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400683", 25, $null)
    ,@("2025-09-26 01:42:45", "金融関連プロジェクトの要件定義や基本設計ができる方(1人月、長期継続案件)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400681", 25, $null)
    ,@("2025-09-26 01:42:45", "【急募】SOLIDWORKS2024での機械設計と製図依頼", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400338", 25, $null)
    ,@("2025-09-26 01:42:45", "当社CTソフトへの機能追加:気道抽出", "システム開発", "3,000,000 円 ~ 5,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400094", 25, $null)
    ,@("2025-09-26 01:42:45", "【SES案件多数/リモート可】フルスタックエンジニア募集(フロント〜バック〜クラウドまで歓迎)", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399877", 25, $null)
    ,@("2025-09-26 01:42:45", "【SES案件多数/リモート可】インフラエンジニア募集(AWS/Linux/NW設計・構築 等歓迎)", "システム開発", "500,000 円 ~ 1,000,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399876", 25, $null)
    ,@("2025-09-26 01:42:45", "【急募】Nuxt3でのWEBページ表示速度改善依頼", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400231", 18, $null)
    ,@("2025-09-26 01:42:45", "限定公開 PR 限定公開の仕事", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5399347", 13, $null)
    ,@("2025-09-26 01:42:45", "【急募】ライフプランシミュレーターのバグ確認と使用感調査", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400626", 10, $null)
    ,@("2025-09-26 01:42:45", "【SalesIQ活用】CRMと連携したリード獲得方法を教えてください", "システム開発", "~ 5,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5400402", 10, $null)
)

# Clear any stale hyperlinks before rebuilding F2:F20 so relationship ids stay in sync
# with the row order below (Range(...).Hyperlinks.Delete() clears the whole sheet in this host).
$ws.Range("F2").Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
    if ($row[7] -ne $null) {
        $ws.Range("H$r").Value = $row[7]
    } else {
        $ws.Range("H$r").Value = ""
    }
    $ws.Hyperlinks.Add($ws.Range("F$r"), $row[5])
    $ws.Range("F$r").Style = "Hyperlink"
    $r++
}
